$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.909.90'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.875.11'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'0.7407"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.34%  '
$ws.Range('D6').Value = "'242.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = "'1.000"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.74%  '
$ws.Range('D9').Value = "'0.07224"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('E10').Value = '  -3.94%  '
$ws.Range('D11').Value = "'0.08341"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.22%  '
$ws.Range('D12').Value = "'0.7500"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('D13').Value = '1.888.57'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').Value = "'5.390"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').Value = "'92.30"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').Value = "'6.115"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.62%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '29.904.76'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = "'246.73"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('D19').Value = "'13.58"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.52%  '
$ws.Range('D20').Value = "'0.000007843"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').Value = "'1.001"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('D22').Value = '2.141.12'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').Value = "'8.000"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').Value = "'0.9991"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('E25').Value = '  -5.54%  '
$ws.Range('D26').Value = "'9.295"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('D27').Value = "'165.39"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('D28').Value = "'18.66"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').Value = "'4.608"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.60%  '
$ws.Range('D32').Value = "'1.537"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').Value = "'4.230"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.35%  '
$ws.Range('D34').Value = "'0.05345"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('D36').Value = "'0.7496"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = "'2.700"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').Value = "'0.01959"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').Value = "'2.752"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('D41').Value = "'0.4521"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('D42').Value = '1.114.22'
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('D43').Value = "'6.142"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.23%  '
$ws.Range('D44').Value = "'72.50"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('D45').Value = "'0.8622"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('D46').Value = "'104.29"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = "'1.865"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').Value = "'7.611"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').Value = "'9.512"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('D51').Value = '2.037.32'
$ws.Range('E51').Value = '  -0.45%  '
